$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("met_metadata")
$wsBatch = $wb.Worksheets.Item("batch_example")

# Update the collision energy values on batch_example sheet to a text value "20 eV"
$wsBatch.Range("C2").Value = "20 eV"
$wsBatch.Range("C3").Value = "20 eV"

# Change selections
$wsMeta.Range("C1").Select() | Out-Null
$wsBatch.Range("C4").Select() | Out-Null

# Activate batch_example sheet (becomes the active tab)
$wsBatch.Activate()
